$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (row 33: activity, row 34: weekly_activity, both 2025-02-11) ---
$ws.Range("F33").Value = $true
$ws.Range("K33").Value = $true
$ws.Range("F34").Value = $true

# --- Append new rows for 2025-02-12 ---
# Format the new Date cells as text first so the "YYYY-MM-DD" strings are
# stored as literal text (matching the rest of column A) instead of being
# auto-converted to date serial numbers.
$ws.Range("A35:A37").NumberFormat = "@"

# Row 35: sleep
$ws.Range("A35").Value = "2025-02-12"
$ws.Range("B35").Value = "sleep"
$ws.Range("C35").Value = $false
$ws.Range("D35").Value = $false
$ws.Range("E35").Value = $true
$ws.Range("F35").Value = $false
$ws.Range("G35").Value = $false
$ws.Range("H35").Value = $true
$ws.Range("I35").Value = $true
$ws.Range("J35").Value = $false
$ws.Range("K35").Value = $true
$ws.Range("L35").Value = $true
$ws.Range("M35").Value = $true
$ws.Range("N35").Value = $true
$ws.Range("O35").Value = $true

# Row 36: activity
$ws.Range("A36").Value = "2025-02-12"
$ws.Range("B36").Value = "activity"
$ws.Range("C36").Value = $false
$ws.Range("D36").Value = $false
$ws.Range("E36").Value = $false
$ws.Range("F36").Value = $false
$ws.Range("G36").Value = $false
$ws.Range("H36").Value = $true
$ws.Range("I36").Value = $true
$ws.Range("J36").Value = $true
$ws.Range("K36").Value = $false
$ws.Range("L36").Value = $true
$ws.Range("M36").Value = $true
$ws.Range("N36").Value = $false
$ws.Range("O36").Value = $false

# Row 37: weekly_activity
$ws.Range("A37").Value = "2025-02-12"
$ws.Range("B37").Value = "weekly_activity"
$ws.Range("C37").Value = $false
$ws.Range("D37").Value = $false
$ws.Range("E37").Value = $true
$ws.Range("F37").Value = $true
$ws.Range("G37").Value = $false
$ws.Range("H37").Value = $false
$ws.Range("I37").Value = $true
$ws.Range("J37").Value = $true
$ws.Range("K37").Value = $true
$ws.Range("L37").Value = $false
$ws.Range("M37").Value = $true
$ws.Range("N37").Value = $false
$ws.Range("O37").Value = $false
